$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the merged "A:B" column-width definition into separate per-column entries
# (A keeps its existing 30.7109375 width; B/C already have their own correct entries).
$ws.Columns.Item(1).Hidden = $false

# --- Row 1 ---
$ws.Cells.Item(1,1).Clear()
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(1,2).PasteSpecial(-4122)
$ws.Cells.Item(1,2).Value = 'Ementa atual:'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(1,3).PasteSpecial(-4122)
$ws.Cells.Item(1,3).Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Rows.Item(1).AutoFit()

# --- Row 2 ---
$ws.Cells.Item(2,1).Clear()
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(2,2).PasteSpecial(-4122)
$ws.Cells.Item(2,2).Value = 'LOQ4055'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(2,3).PasteSpecial(-4122)
$ws.Cells.Item(2,3).Value = 'LOQ4055'
$ws.Rows.Item(2).AutoFit()

# --- Row 3 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(3,1).PasteSpecial(-4122)
$ws.Cells.Item(3,1).Value = 'Nome:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(3,2).PasteSpecial(-4122)
$ws.Cells.Item(3,2).Value = ' Quimica Inorgânica'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(3,3).PasteSpecial(-4122)
$ws.Cells.Item(3,3).Value = ' Quimica Inorgânica'
$ws.Rows.Item(3).AutoFit()

# --- Row 4 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(4,1).PasteSpecial(-4122)
$ws.Cells.Item(4,1).Value = 'Name:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(4,2).PasteSpecial(-4122)
$ws.Cells.Item(4,2).Value = 'Inorganic Chemistry'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(4,3).PasteSpecial(-4122)
$ws.Cells.Item(4,3).Value = 'Inorganic Chemistry'
$ws.Rows.Item(4).AutoFit()

# --- Row 5 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(5,1).PasteSpecial(-4122)
$ws.Cells.Item(5,1).Value = 'Créditos-aula:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(5,2).PasteSpecial(-4122)
$ws.Cells.Item(5,2).Value = '2'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(5,3).PasteSpecial(-4122)
$ws.Cells.Item(5,3).Value = '2'
$ws.Rows.Item(5).AutoFit()

# --- Row 6 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(6,1).PasteSpecial(-4122)
$ws.Cells.Item(6,1).Value = 'Créditos-trabalho'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(6,2).PasteSpecial(-4122)
$ws.Cells.Item(6,2).Value = '0'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(6,3).PasteSpecial(-4122)
$ws.Cells.Item(6,3).Value = '0'
$ws.Rows.Item(6).AutoFit()

# --- Row 7 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(7,1).PasteSpecial(-4122)
$ws.Cells.Item(7,1).Value = 'Carga horária:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(7,2).PasteSpecial(-4122)
$ws.Cells.Item(7,2).Value = '30 h'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(7,3).PasteSpecial(-4122)
$ws.Cells.Item(7,3).Value = '30 h'
$ws.Rows.Item(7).AutoFit()

# --- Row 8 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(8,1).PasteSpecial(-4122)
$ws.Cells.Item(8,1).Value = 'Ativação:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(8,2).PasteSpecial(-4122)
$ws.Cells.Item(8,2).Value = '01/01/2022'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(8,3).PasteSpecial(-4122)
$ws.Cells.Item(8,3).Value = '01/01/2022'
$ws.Rows.Item(8).AutoFit()

# --- Row 9 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(9,1).PasteSpecial(-4122)
$ws.Cells.Item(9,1).Value = 'Semestre ideal:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(9,2).PasteSpecial(-4122)
$ws.Cells.Item(9,2).Value = 'EQD-3,EQN-4'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(9,3).PasteSpecial(-4122)
$ws.Cells.Item(9,3).Value = 'EQD-3,EQN-4'
$ws.Rows.Item(9).AutoFit()

# --- Row 10 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(10,1).PasteSpecial(-4122)
$ws.Cells.Item(10,1).Value = 'Objetivos:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(10,2).PasteSpecial(-4122)
$ws.Cells.Item(10,2).Value = 'Fornecer aos alunos conceitos fundamentais para compreensão da Química Inorgânica por meio da experimentação, desenvolvendo a capacidade de realizarem práticas no laboratório que estimulem o seu pensamento científico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de caráter inorgânico com interesse industrial.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(10,3).PasteSpecial(-4122)
$ws.Cells.Item(10,3).Value = 'Fornecer aos alunos conceitos fundamentais para compreensão da Química Inorgânica por meio da experimentação, desenvolvendo a capacidade de realizarem práticas no laboratório que estimulem o seu pensamento científico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de caráter inorgânico com interesse industrial.'
$ws.Rows.Item(10).RowHeight = 60

# --- Row 11 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(11,1).PasteSpecial(-4122)
$ws.Cells.Item(11,1).Value = 'Objectives:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(11,2).PasteSpecial(-4122)
$ws.Cells.Item(11,2).Value = 'Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(11,3).PasteSpecial(-4122)
$ws.Cells.Item(11,3).Value = 'Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.'
$ws.Rows.Item(11).RowHeight = 60

# --- Row 12 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(12,1).PasteSpecial(-4122)
$ws.Cells.Item(12,1).Value = 'Docentes responsáveis:'
$ws.Cells.Item(12,2).Clear()
$ws.Cells.Item(12,3).Clear()
$ws.Rows.Item(12).AutoFit()

# --- Row 13 ---
$ws.Cells.Item(13,1).Clear()
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(13,2).PasteSpecial(-4122)
$ws.Cells.Item(13,2).Value = '5840705 - Maria Lúcia Caetano Pinto da Silva'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(13,3).PasteSpecial(-4122)
$ws.Cells.Item(13,3).Value = '5840705 - Maria Lúcia Caetano Pinto da Silva'
$ws.Rows.Item(13).AutoFit()

# --- Row 14 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(14,1).PasteSpecial(-4122)
$ws.Cells.Item(14,1).Value = 'Programa resumido:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(14,2).PasteSpecial(-4122)
$ws.Cells.Item(14,2).Value = 'Compostos de Coordenação. Materiais inorgânicos de interesse industrial. Purificação e Identificação de Compostos Inorgânicos. Síntese de sais e obtenção de Compostos de Alumínio.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(14,3).PasteSpecial(-4122)
$ws.Cells.Item(14,3).Value = 'Compostos de Coordenação. Materiais inorgânicos de interesse industrial. Purificação e Identificação de Compostos Inorgânicos. Síntese de sais e obtenção de Compostos de Alumínio.'
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(15,1).PasteSpecial(-4122)
$ws.Cells.Item(15,1).Value = 'Short syllabus:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(15,2).PasteSpecial(-4122)
$ws.Cells.Item(15,2).Value = 'Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(15,3).PasteSpecial(-4122)
$ws.Cells.Item(15,3).Value = 'Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value = 'Programa:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(16,2).PasteSpecial(-4122)
$ws.Cells.Item(16,2).Value = 'Compostos de Coordenação: Estrutura, ligações, reações e aplicações. Exemplos e aplicações de materiais inorgânicos de interesse industrial. Sínteses: Sal Simples, Sal Duplo e Sal Complexo. Preparação de Compostos de Alumínio.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(16,3).PasteSpecial(-4122)
$ws.Cells.Item(16,3).Value = 'Compostos de Coordenação: Estrutura, ligações, reações e aplicações. Exemplos e aplicações de materiais inorgânicos de interesse industrial. Sínteses: Sal Simples, Sal Duplo e Sal Complexo. Preparação de Compostos de Alumínio.'
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(17,1).PasteSpecial(-4122)
$ws.Cells.Item(17,1).Value = 'Syllabus:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(17,2).PasteSpecial(-4122)
$ws.Cells.Item(17,2).Value = 'Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(17,3).PasteSpecial(-4122)
$ws.Cells.Item(17,3).Value = 'Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.'
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(18,1).PasteSpecial(-4122)
$ws.Cells.Item(18,1).Value = 'Avaliação:'
$ws.Cells.Item(18,2).Clear()
$ws.Cells.Item(18,3).Clear()
$ws.Rows.Item(18).AutoFit()

# --- Row 19 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(19,1).PasteSpecial(-4122)
$ws.Cells.Item(19,1).Value = 'Método:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(19,2).PasteSpecial(-4122)
$ws.Cells.Item(19,2).Value = 'Serão oferecidas aulas expositivas e práticas.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(19,3).PasteSpecial(-4122)
$ws.Cells.Item(19,3).Value = 'Serão oferecidas aulas expositivas e práticas.'
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(20,1).PasteSpecial(-4122)
$ws.Cells.Item(20,1).Value = 'Critério:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(20,2).PasteSpecial(-4122)
$ws.Cells.Item(20,2).Value = 'Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(20,3).PasteSpecial(-4122)
$ws.Cells.Item(20,3).Value = 'Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita.'
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(21,1).PasteSpecial(-4122)
$ws.Cells.Item(21,1).Value = 'Norma de recuperação:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(21,2).PasteSpecial(-4122)
$ws.Cells.Item(21,2).Value = 'Será realizada uma prova escrita envolvendo o conteúdo do semestre todo.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(21,3).PasteSpecial(-4122)
$ws.Cells.Item(21,3).Value = 'Será realizada uma prova escrita envolvendo o conteúdo do semestre todo.'
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(22,1).PasteSpecial(-4122)
$ws.Cells.Item(22,1).Value = 'Bibliografia:'
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(22,2).PasteSpecial(-4122)
$ws.Cells.Item(22,2).Value = 'CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981.LEE, J. D., tradução Química Inorgânica não tão concisa da 5ª edição inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Química Inorgânica tradução da 4ª edição. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Química - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3ª ed., 1973.'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(22,3).PasteSpecial(-4122)
$ws.Cells.Item(22,3).Value = 'CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981.LEE, J. D., tradução Química Inorgânica não tão concisa da 5ª edição inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Química Inorgânica tradução da 4ª edição. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Química - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3ª ed., 1973.'
$ws.Rows.Item(22).RowHeight = 120

# --- Row 23 ---
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(23,1).PasteSpecial(-4122)
$ws.Cells.Item(23,1).Value = 'Requisitos:'
$ws.Cells.Item(23,2).Clear()
$ws.Cells.Item(23,3).Clear()
$ws.Rows.Item(23).AutoFit()

# --- Row 24 ---
$ws.Cells.Item(24,1).Clear()
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(24,2).PasteSpecial(-4122)
$ws.Cells.Item(24,2).Value = 'LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(24,3).PasteSpecial(-4122)
$ws.Cells.Item(24,3).Value = 'LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'
$ws.Rows.Item(24).RowHeight = 30

$excel.CutCopyMode = $false